$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.804.41'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").Value = '2.034.57'
$ws.Range("E3").Value = '  -1.14%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '''227.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").Value = '''60.21'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.08%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -2.07%  '
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").Value = '2.337.14'
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").Value = '''14.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.80%  '
$ws.Range("D14").Value = '''21.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").Value = '''0.759'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").Value = '2.033.47'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = '37.795.90'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '''69.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = '''5.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.48%  '
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("E21").Value = '  -1.42%  '
$ws.Range("D22").Value = '''223.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("D26").Value = '''167.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("D27").Value = '''9.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").Value = '''18.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  -3.86%  '
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").Value = '''2.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.08%  '
$ws.Range("D33").Value = '''4.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.38%  '
$ws.Range("E34").Value = '  -0.17%  '
$ws.Range("D35").Value = '''4.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.26%  '
$ws.Range("D36").Value = '''6.31'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("D37").Value = '''2.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.36%  '
$ws.Range("D38").Value = '''3.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.33%  '
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("D40").Value = '''17.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.77%  '
$ws.Range("D41").Value = '1.533.60'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '''0.0216'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = '''96.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("D45").Value = '''0.0914'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("D46").Value = '''1.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.55%  '
$ws.Range("D47").Value = '''4.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.63%  '
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("D49").Value = '''7.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").Value = '2.225.52'
$ws.Range("E51").Value = '  -1.07%  '
